$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before the existing row 765, shifting all the following
# rows (765-816) down to (772-823).
$ws.Rows("765:771").Insert()

# Populate the 7 newly inserted rows (765-771) with their data.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Tipo,
#          G Producto ID, H Producto, I Categoría ID, J Categoría,
#          K Variedad, L Calidad, M Volumen, N Precio minimo,
#          O Precio maximo, P Precio promedio ponderado,
#          Q Unidad de comercializacion, R Origen, S Precio $/Kg, T Kg/unidad

$ws.Cells.Item(765,1).Value2 = 6
$ws.Cells.Item(765,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(765,3).Value2 = "Metropolitana"
$ws.Cells.Item(765,4).Value2 = 44516
$ws.Cells.Item(765,5).Value2 = 13
$ws.Cells.Item(765,6).Value2 = "Fruta"
$ws.Cells.Item(765,7).Value2 = 100104
$ws.Cells.Item(765,8).Value2 = "Frutos de pepita"
$ws.Cells.Item(765,9).Value2 = 100104005
$ws.Cells.Item(765,10).Value2 = "Pera"
$ws.Cells.Item(765,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(765,12).Value2 = "Especial"
$ws.Cells.Item(765,13).Value2 = 10
$ws.Cells.Item(765,14).Value2 = 240000
$ws.Cells.Item(765,15).Value2 = 240000
$ws.Cells.Item(765,16).Value2 = 240000
$ws.Cells.Item(765,17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(765,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(765,19).Value2 = 533
$ws.Cells.Item(765,20).Value2 = 450

$ws.Cells.Item(766,1).Value2 = 6
$ws.Cells.Item(766,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(766,3).Value2 = "Metropolitana"
$ws.Cells.Item(766,4).Value2 = 44516
$ws.Cells.Item(766,5).Value2 = 13
$ws.Cells.Item(766,6).Value2 = "Fruta"
$ws.Cells.Item(766,7).Value2 = 100104
$ws.Cells.Item(766,8).Value2 = "Frutos de pepita"
$ws.Cells.Item(766,9).Value2 = 100104005
$ws.Cells.Item(766,10).Value2 = "Pera"
$ws.Cells.Item(766,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(766,12).Value2 = "Primera"
$ws.Cells.Item(766,13).Value2 = 17
$ws.Cells.Item(766,14).Value2 = 220000
$ws.Cells.Item(766,15).Value2 = 220000
$ws.Cells.Item(766,16).Value2 = 220000
$ws.Cells.Item(766,17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(766,18).Value2 = "Paine"
$ws.Cells.Item(766,19).Value2 = 489
$ws.Cells.Item(766,20).Value2 = 450

$ws.Cells.Item(767,1).Value2 = 6
$ws.Cells.Item(767,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(767,3).Value2 = "Metropolitana"
$ws.Cells.Item(767,4).Value2 = 44516
$ws.Cells.Item(767,5).Value2 = 13
$ws.Cells.Item(767,6).Value2 = "Fruta"
$ws.Cells.Item(767,7).Value2 = 100104
$ws.Cells.Item(767,8).Value2 = "Frutos de pepita"
$ws.Cells.Item(767,9).Value2 = 100104005
$ws.Cells.Item(767,10).Value2 = "Pera"
$ws.Cells.Item(767,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(767,12).Value2 = "Primera"
$ws.Cells.Item(767,13).Value2 = 14
$ws.Cells.Item(767,14).Value2 = 220000
$ws.Cells.Item(767,15).Value2 = 220000
$ws.Cells.Item(767,16).Value2 = 220000
$ws.Cells.Item(767,17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(767,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(767,19).Value2 = 489
$ws.Cells.Item(767,20).Value2 = 450

$ws.Cells.Item(768,1).Value2 = 6
$ws.Cells.Item(768,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(768,3).Value2 = "Metropolitana"
$ws.Cells.Item(768,4).Value2 = 44516
$ws.Cells.Item(768,5).Value2 = 13
$ws.Cells.Item(768,6).Value2 = "Fruta"
$ws.Cells.Item(768,7).Value2 = 100104
$ws.Cells.Item(768,8).Value2 = "Frutos de pepita"
$ws.Cells.Item(768,9).Value2 = 100104005
$ws.Cells.Item(768,10).Value2 = "Pera"
$ws.Cells.Item(768,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(768,12).Value2 = "Segunda"
$ws.Cells.Item(768,13).Value2 = 15
$ws.Cells.Item(768,14).Value2 = 180000
$ws.Cells.Item(768,15).Value2 = 180000
$ws.Cells.Item(768,16).Value2 = 180000
$ws.Cells.Item(768,17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(768,18).Value2 = "Paine"
$ws.Cells.Item(768,19).Value2 = 400
$ws.Cells.Item(768,20).Value2 = 450

$ws.Cells.Item(769,1).Value2 = 6
$ws.Cells.Item(769,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(769,3).Value2 = "Metropolitana"
$ws.Cells.Item(769,4).Value2 = 44516
$ws.Cells.Item(769,5).Value2 = 13
$ws.Cells.Item(769,6).Value2 = "Fruta"
$ws.Cells.Item(769,7).Value2 = 100104
$ws.Cells.Item(769,8).Value2 = "Frutos de pepita"
$ws.Cells.Item(769,9).Value2 = 100104005
$ws.Cells.Item(769,10).Value2 = "Pera"
$ws.Cells.Item(769,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(769,12).Value2 = "Segunda"
$ws.Cells.Item(769,13).Value2 = 18
$ws.Cells.Item(769,14).Value2 = 200000
$ws.Cells.Item(769,15).Value2 = 200000
$ws.Cells.Item(769,16).Value2 = 200000
$ws.Cells.Item(769,17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(769,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(769,19).Value2 = 444
$ws.Cells.Item(769,20).Value2 = 450

$ws.Cells.Item(770,1).Value2 = 6
$ws.Cells.Item(770,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(770,3).Value2 = "Metropolitana"
$ws.Cells.Item(770,4).Value2 = 44516
$ws.Cells.Item(770,5).Value2 = 13
$ws.Cells.Item(770,6).Value2 = "Fruta"
$ws.Cells.Item(770,7).Value2 = 100104
$ws.Cells.Item(770,8).Value2 = "Frutos de pepita"
$ws.Cells.Item(770,9).Value2 = 100104005
$ws.Cells.Item(770,10).Value2 = "Pera"
$ws.Cells.Item(770,11).Value2 = "Winter Nelis"
$ws.Cells.Item(770,12).Value2 = "Especial"
$ws.Cells.Item(770,13).Value2 = 5
$ws.Cells.Item(770,14).Value2 = 300000
$ws.Cells.Item(770,15).Value2 = 300000
$ws.Cells.Item(770,16).Value2 = 300000
$ws.Cells.Item(770,17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(770,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(770,19).Value2 = 667
$ws.Cells.Item(770,20).Value2 = 450

$ws.Cells.Item(771,1).Value2 = 6
$ws.Cells.Item(771,2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(771,3).Value2 = "Metropolitana"
$ws.Cells.Item(771,4).Value2 = 44516
$ws.Cells.Item(771,5).Value2 = 13
$ws.Cells.Item(771,6).Value2 = "Fruta"
$ws.Cells.Item(771,7).Value2 = 100104
$ws.Cells.Item(771,8).Value2 = "Frutos de pepita"
$ws.Cells.Item(771,9).Value2 = 100104005
$ws.Cells.Item(771,10).Value2 = "Pera"
$ws.Cells.Item(771,11).Value2 = "Winter Nelis"
$ws.Cells.Item(771,12).Value2 = "Segunda"
$ws.Cells.Item(771,13).Value2 = 8
$ws.Cells.Item(771,14).Value2 = 220000
$ws.Cells.Item(771,15).Value2 = 220000
$ws.Cells.Item(771,16).Value2 = 220000
$ws.Cells.Item(771,17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(771,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(771,19).Value2 = 489
$ws.Cells.Item(771,20).Value2 = 450
